# Add the new "2022" data column (column J) to the poverty-rate table,
# mirroring the formatting already used by the neighbouring "2021" column
# (column I), and leave the sheet with D1 selected — matching the
# published edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats / xlVAlignBottom constants (PowerShell/VBA numeric values,
# used directly since this host doesn't expose the Excel enum names).
$xlPasteFormats = -4122
$xlVAlignBottom = -4107

function Set-LikeNeighbor {
    param(
        [string]$SourceCell,
        [string]$TargetCell,
        $Value
    )

    $ws.Range($SourceCell).Copy() | Out-Null
    $ws.Range($TargetCell).PasteSpecial($xlPasteFormats) | Out-Null
    # The source column's cells are vertically centred; the new column's
    # cells are not, so drop the inherited vertical centring.
    $ws.Range($TargetCell).VerticalAlignment = $xlVAlignBottom

    if ($null -ne $Value) {
        $ws.Range($TargetCell).Value = $Value
    }
}

$ws.Application.CutCopyMode = $false

# Header: J4 = 2022, formatted exactly like the other year headers (D4:I4).
Set-LikeNeighbor -SourceCell "I4" -TargetCell "J4" -Value 2022

# Row 5 ("Kyrgyz Republic" total line) uses the bold header-ish row style.
Set-LikeNeighbor -SourceCell "I5" -TargetCell "J5" -Value 47.4

# Plain (regular) data / blank sub-heading rows.
$plainRows = @(6, 7, 9, 11, 13, 14, 15, 16, 17, 19, 21, 24)
$plainValues = @{
    6  = $null
    7  = 47.9
    9  = $null
    11 = 41.4
    13 = $null
    14 = 39.5
    15 = 51.9
    16 = $null
    17 = 69.900000000000006
    19 = 42.5
    21 = 42.5
    24 = 38.9
}
foreach ($r in $plainRows) {
    Set-LikeNeighbor -SourceCell "I6" -TargetCell ("J" + $r) -Value $plainValues[$r]
}

# Rows whose "2021" figure is styled with the 0.0 number format (bold
# sub-totals); match that for the new "2022" figure too.
$numfmtRows = @(8, 10, 12, 18, 20, 22, 23)
$numfmtValues = @{
    8  = 46.9
    10 = 56.8
    12 = 39
    18 = 61
    20 = 54
    22 = 45.8
    23 = 38.1
}
foreach ($r in $numfmtRows) {
    Set-LikeNeighbor -SourceCell "I8" -TargetCell ("J" + $r) -Value $numfmtValues[$r]
}

# Last table row (25) carries the bottom border.
Set-LikeNeighbor -SourceCell "I25" -TargetCell "J25" -Value 38.700000000000003

$ws.Application.CutCopyMode = $false

# Leave the selection on D1, as in the published workbook.
$ws.Range("D1").Select() | Out-Null
